# Auto-generated Excel COM-interop script applying the Gungnir_Profits.xlsx diff.
# For each affected (sheet, row) we set the changed H-N columns to their new values,
# and clear any cell that the diff removes entirely (M7/N7 in CUL, N122 in ARM).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2796.074
$ws.Range("I113").Value = 2045
$ws.Range("J113").Value = 3312.4375
$ws.Range("K113").Value = 2045
$ws.Range("L113").Value = 3312.4375
$ws.Range("M113").Value = 1209
$ws.Range("N113").Value = -9820.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4769.31
$ws.Range("I32").Value = 4011.2908
$ws.Range("J32").Value = 9425.714
$ws.Range("K32").Value = 4011.2908
$ws.Range("L32").Value = 9425.714
$ws.Range("M32").Value = -3724.2908
$ws.Range("N32").Value = -9999.714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9288169
$ws.Range("I45").Value = 13374248
$ws.Range("K45").Value = 13374248
$ws.Range("M45").Value = -13373871

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 19800
$ws.Range("J64").Value = 19800
$ws.Range("L64").Value = 19800
$ws.Range("N64").Value = -20296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 19800
$ws.Range("J67").Value = 19800
$ws.Range("L67").Value = 19800
$ws.Range("N67").Value = -21516

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 13890265
$ws.Range("I74").Value = 16668126
$ws.Range("J74").Value = 957
$ws.Range("K74").Value = 16668126
$ws.Range("L74").Value = 957
$ws.Range("M74").Value = -16667252
$ws.Range("N74").Value = -2705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 13890265
$ws.Range("I77").Value = 16668126
$ws.Range("J77").Value = 957
$ws.Range("K77").Value = 83340630
$ws.Range("L77").Value = 4785
$ws.Range("M77").Value = -83336262
$ws.Range("N77").Value = -13521

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1498.1818
$ws.Range("I122").Value = 1498.1818
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4494.5454
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2044.5454
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 20411556
$ws.Range("I132").Value = 29412904
$ws.Range("K132").Value = 88238712
$ws.Range("M132").Value = -88236182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3015.48
$ws.Range("I31").Value = 2157.092
$ws.Range("J31").Value = 5733.7085
$ws.Range("K31").Value = 2157.092
$ws.Range("L31").Value = 5733.7085
$ws.Range("M31").Value = -1862.092
$ws.Range("N31").Value = -6323.7085

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3015.48
$ws.Range("I34").Value = 2157.092
$ws.Range("J34").Value = 5733.7085
$ws.Range("K34").Value = 2157.092
$ws.Range("L34").Value = 5733.7085
$ws.Range("M34").Value = -1955.092
$ws.Range("N34").Value = -6137.7085

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 5237.3335
$ws.Range("I105").Value = 5885.0527
$ws.Range("K105").Value = 5885.0527
$ws.Range("M105").Value = -4138.0527

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 8334468.5
$ws.Range("I122").Value = 11905452
$ws.Range("J122").Value = 2173.7778
$ws.Range("K122").Value = 35716356
$ws.Range("L122").Value = 6521.3334
$ws.Range("M122").Value = -35713906
$ws.Range("N122").Value = -11421.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7249694.5
$ws.Range("I132").Value = 923.3939
$ws.Range("J132").Value = 25650422
$ws.Range("K132").Value = 2770.1817
$ws.Range("L132").Value = 76951266
$ws.Range("M132").Value = -240.1817000000001
$ws.Range("N132").Value = -76956326

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 986.5
$ws.Range("I34").Value = 335.57144
$ws.Range("J34").Value = 1897.8
$ws.Range("K34").Value = 1006.71432
$ws.Range("L34").Value = 5693.4
$ws.Range("M34").Value = -922.71432
$ws.Range("N34").Value = -5861.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1093.75
$ws.Range("J39").Value = 1800
$ws.Range("L39").Value = 5400
$ws.Range("N39").Value = -5988

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1513.3334
$ws.Range("I55").Value = 866.6667
$ws.Range("J55").Value = 1675
$ws.Range("K55").Value = 2600.0001
$ws.Range("L55").Value = 5025
$ws.Range("M55").Value = -2423.0001
$ws.Range("N55").Value = -5379

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 873.17
$ws.Range("J131").Value = 877.70105
$ws.Range("L131").Value = 2633.10315
$ws.Range("N131").Value = -12713.10315

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 19237282
$ws.Range("I122").Value = 35724124
$ws.Range("K122").Value = 107172372
$ws.Range("M122").Value = -107169922

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 3000
$ws.Range("J3").Value = 3000
$ws.Range("L3").Value = 3000
$ws.Range("N3").Value = -3224

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H15").Value = 3000
$ws.Range("J15").Value = 3000
$ws.Range("L15").Value = 3000
$ws.Range("N15").Value = -3340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6506.089
$ws.Range("I122").Value = 7790.273
$ws.Range("J122").Value = 2974.5833
$ws.Range("K122").Value = 23370.819
$ws.Range("L122").Value = 8923.749899999999
$ws.Range("M122").Value = -20920.819
$ws.Range("N122").Value = -13823.7499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 13518832
$ws.Range("I132").Value = 23257196
$ws.Range("J132").Value = 10777.839
$ws.Range("K132").Value = 69771588
$ws.Range("L132").Value = 32333.517
$ws.Range("M132").Value = -69769058
$ws.Range("N132").Value = -37393.517

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2744.2856
$ws.Range("I122").Value = 2250.4443
$ws.Range("J122").Value = 3114.6667
$ws.Range("K122").Value = 6751.3329
$ws.Range("L122").Value = 9344.000100000001
$ws.Range("M122").Value = -4301.3329
$ws.Range("N122").Value = -14244.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5896737.5
$ws.Range("I132").Value = 15703.328
$ws.Range("J132").Value = 33341564
$ws.Range("K132").Value = 47109.984
$ws.Range("L132").Value = 100024692
$ws.Range("M132").Value = -44579.984
$ws.Range("N132").Value = -100029752

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2155.0603
$ws.Range("I136").Value = 2331.8103
$ws.Range("J136").Value = 1745
$ws.Range("K136").Value = 6995.4309
$ws.Range("L136").Value = 5235
$ws.Range("M136").Value = -4445.4309
$ws.Range("N136").Value = -10335
